{"js": "// Apply \"hybrid bold + color\" highlighting to quantitative impact metrics\n// (percentages, dollar amounts, large numbers) across the resume body.\n//\n// Strategy: for each target paragraph (matched by its current, unedited\n// full text) find each metric substring with Range.search() and apply\n// bold + the accent color (#2C3E50) to just that sub-range. Word/Office.js\n// automatically splits the run and preserves the surrounding plain-text\n// runs (including the xml:space=\"preserve\" needed for the space-only\n// segments), which is exactly the run structure produced by the diff.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Map of \"paragraph's current full text\" -> ordered list of metric\n// substrings that must become bold + colored within that paragraph.\nconst EDITS = [\n  {\n    match: \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    metrics: [\"23%\", \"64%\"],\n  },\n  {\n    match: \"\u2022 Utilized advanced sampling methods to decrease survey margin of error from \u00b14.2% to \u00b12.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes\",\n    metrics: [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"],\n  },\n  {\n    match: \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n    metrics: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    match: \"\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion\",\n    metrics: [\"$2\"],\n  },\n  {\n    match: \"\u2022 Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from \u00b14.2% to \u00b12.1%\",\n    metrics: [\"\u00b14.2%\", \"\u00b12.1%\"],\n  },\n  {\n    match: \"\u2022 Increased voter turnout prediction accuracy from 71% to 87%\",\n    metrics: [\"71%\", \"87%\"],\n  },\n  {\n    match: \"\u2022 Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%\",\n    metrics: [\"34%\", \"28%\"],\n  },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nfor (const edit of EDITS) {\n  const paragraph = paragraphs.items.find((p) => p.text === edit.match);\n  if (!paragraph) {\n    throw new Error(\"Could not locate paragraph: \" + edit.match);\n  }\n\n  for (const metric of edit.metrics) {\n    const found = paragraph.search(metric, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n\n    if (found.items.length === 0) {\n      throw new Error(\"Could not locate metric '\" + metric + \"' in paragraph: \" + edit.match);\n    }\n\n    const range = found.items[0];\n    range.font.bold = true;\n    range.font.color = HIGHLIGHT_COLOR;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply \"hybrid bold + color\" highlighting to quantitative impact metrics\n# (percentages, dollar amounts, large numbers) across the resume body.\n#\n# Strategy: for each target paragraph (identified by its fixed position in\n# the document, verified against the expected leading text) find each\n# metric substring with Range.Find.Execute() scoped to that paragraph, and\n# apply Bold + the accent color (#2C3E50) to just that sub-range. Word\n# automatically splits the run and preserves the surrounding plain-text\n# runs (including the space-only segments), producing the same run\n# structure as the target diff.\n\n$d = $word.ActiveDocument\n\n# #2C3E50 as a WdColor (0xBBGGRR) value for Font.Color.\n$HighlightColor = 5258796\n\nfunction Set-MetricHighlight($paragraphIndex, $expectedPrefix, $metric) {\n    $para = $d.Paragraphs.Item($paragraphIndex)\n    $prefixRange = $para.Range.Duplicate\n    $prefixLen = [Math]::Min($expectedPrefix.Length, $prefixRange.Text.Length)\n    $prefixRange.End = $prefixRange.Start + $prefixLen\n    if ($prefixRange.Text -ne $expectedPrefix) {\n        throw \"Paragraph $paragraphIndex does not start with expected text '$expectedPrefix' (got '$($prefixRange.Text)')\"\n    }\n\n    $rng = $para.Range\n    $found = $rng.Find.Execute($metric)\n    if (-not $found) {\n        throw \"Could not find metric '$metric' in paragraph $paragraphIndex\"\n    }\n    $rng.Font.Bold = $true\n    $rng.Font.Color = $HighlightColor\n}\n\n# Paragraph 10: \"\u2022 Discovered systematic race coding errors ... from 23% to 64%\"\nSet-MetricHighlight 10 \"\u2022 Discovered systematic race coding errors\" \"23%\"\nSet-MetricHighlight 10 \"\u2022 Discovered systematic race coding errors\" \"64%\"\n\n# Paragraph 12: \"\u2022 Utilized advanced sampling methods ... \u00b14.2% to \u00b12.1% ... 71% to 87% ...\"\nSet-MetricHighlight 12 \"\u2022 Utilized advanced sampling methods\" \"\u00b14.2%\"\nSet-MetricHighlight 12 \"\u2022 Utilized advanced sampling methods\" \"\u00b12.1%\"\nSet-MetricHighlight 12 \"\u2022 Utilized advanced sampling methods\" \"71%\"\nSet-MetricHighlight 12 \"\u2022 Utilized advanced sampling methods\" \"87%\"\n\n# Paragraph 13: \"\u2022 Trigonometric algorithm for boundary estimation ... 73.5% ... $4.7M ...\"\nSet-MetricHighlight 13 \"\u2022 Trigonometric algorithm for boundary estimation\" \"73.5%\"\nSet-MetricHighlight 13 \"\u2022 Trigonometric algorithm for boundary estimation\" \"$4.7M\"\n\n# Paragraph 14: \"\u2022 Built real-time FEC analysis systems ... valued over $2 trillion\"\nSet-MetricHighlight 14 \"\u2022 Built real-time FEC analysis systems\" \"$2\"\n\n# Paragraph 50: \"\u2022 Predictive excellence: ... \u00b14.2% to \u00b12.1%\"\nSet-MetricHighlight 50 \"\u2022 Predictive excellence\" \"\u00b14.2%\"\nSet-MetricHighlight 50 \"\u2022 Predictive excellence\" \"\u00b12.1%\"\n\n# Paragraph 51: \"\u2022 Increased voter turnout prediction accuracy from 71% to 87%\"\nSet-MetricHighlight 51 \"\u2022 Increased voter turnout prediction accuracy\" \"71%\"\nSet-MetricHighlight 51 \"\u2022 Increased voter turnout prediction accuracy\" \"87%\"\n\n# Paragraph 53: \"\u2022 Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%\"\nSet-MetricHighlight 53 \"\u2022 Methodological advancement\" \"34%\"\nSet-MetricHighlight 53 \"\u2022 Methodological advancement\" \"28%\"\n"}
